$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Two weeks' worth of rows (old A14:B14 = 45361.../1080 and old A15:B15 =
# 45368.../340) were removed entirely; everything below shifts up by two.
$ws1.Rows.Item(14).Delete()
$ws1.Rows.Item(14).Delete()

# The row that used to be row 16 (45375.../810) is now row 14; its
# requested quantity changes from 810 down to 50.
$ws1.Cells.Item(14, 2).Value = 50

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Cells.Item(9, 2).Value = 50
